$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with newly scraped values.
# Every cell in this sheet is stored as text (coinranking.com formats the
# price with dotted thousands separators, e.g. "64.021.35"). Some of the new
# Price strings happen to look like plain decimals (e.g. "597.20"); setting
# .Value on those directly would have Excel interpret them as numbers and
# silently drop the trailing zero / change the cell type, so those specific
# cells are pre-formatted as Text to keep the exact original string.

$ws.Range('D2').Value = '64.021.35'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '3.095.39'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.20'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.20'
$ws.Range('E6').Value = '  +1.57%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('D9').Value = '3.095.16'
$ws.Range('E9').Value = '  -1.89%  '
$ws.Range('E10').Value = '  -2.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.97'
$ws.Range('E11').Value = '  -0.96%  '
$ws.Range('E12').Value = '  -3.21%  '
$ws.Range('E13').Value = '  -4.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.91'
$ws.Range('E14').Value = '  -4.70%  '
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('D16').Value = '3.606.52'
$ws.Range('E16').Value = '  -1.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.22'
$ws.Range('E17').Value = '  -1.40%  '
$ws.Range('D18').Value = '63.911.65'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').Value = '3.094.22'
$ws.Range('E19').Value = '  -1.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '480.93'
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.49'
$ws.Range('E21').Value = '  -3.85%  '
$ws.Range('E22').Value = '  -5.01%  '
$ws.Range('E23').Value = '  -1.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.50'
$ws.Range('E24').Value = '  +3.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.65'
$ws.Range('E25').Value = '  -1.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.90'
$ws.Range('E26').Value = '  -5.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.75'
$ws.Range('E27').Value = '  +8.50%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.65'
$ws.Range('E29').Value = '  +2.42%  '
$ws.Range('E30').Value = '  -1.58%  '
$ws.Range('E31').Value = '  -0.53%  '
$ws.Range('E32').Value = '  -1.69%  '
$ws.Range('E33').Value = '  -4.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.21'
$ws.Range('E34').Value = '  -2.31%  '
$ws.Range('E35').Value = '  -4.25%  '
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.30'
$ws.Range('E38').Value = '  -6.70%  '
$ws.Range('E39').Value = '  -3.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '51.03'
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('E41').Value = '  -1.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '444.42'
$ws.Range('E42').Value = '  -5.38%  '
$ws.Range('E43').Value = '  -3.80%  '
$ws.Range('E44').Value = '  -4.63%  '
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('D46').Value = '2.837.99'
$ws.Range('E46').Value = '  -2.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.09'
$ws.Range('E47').Value = '  +3.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.21'
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '26.08'
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('E51').Value = '  -2.58%  '
